$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    "Turkish Radio and Television Corporation (TRT)",
    "TV Azteca",
    "TV Globo",
    "TV Markiza",
    "TV Nova (Czech Republic)",
    "TV Prima",
    "TV Puls",
    "TV Visjon Norge",
    "TV2 Denmark",
    "TV2 Hungary",
    "TV2 Norway",
    "TV3 Group (All Media Baltics)",
    "TV3 Lithuania",
    "TV3 Network Ghana",
    "TV4 (Sweden)",
    "TV5 Monde",
    "TV8 Turkey",
    "TVN Discovery Group",
    "tvOne Indonesia",
    "TVP (Telewizja Polska)",
    "TVRI (Televisi Republik Indonesia)"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
